# "tag delta, move to echo"
#
# Adds the three requirements-summary worksheets (PASS / FAIL / NULL) after
# the existing "provenance" sheet, each carrying a single bold/blue header
# cell, and refreshes the "last generated" timestamp on the provenance sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- provenance sheet: bump the generated-on timestamp (B12) -------------
$ws1.Range("B12").Value = 43431.21433859861

# Same header/footer codes as the "provenance" sheet, split across the
# header/footer section properties (each property setter supplies its own
# leading "&X" section code, so the literal text must omit it).
$headerCenter = "&12&A"
$footerLeft   = "&8&T`n&8&D&C &P / &N&R&8&Z`n&8&F"

# "Section header" formatting (bold blue) already used elsewhere on the
# provenance sheet (e.g. A3 "Workbook created by") - reuse it via a
# formats-only copy/paste instead of rebuilding the font from scratch, so
# no duplicate font/cellXf entries get minted in styles.xml.
$xlPasteFormats = -4122
$styleSource = $ws1.Range("A3")

# --- new sheet: requirements - PASS ---------------------------------------
$wsPass = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$wsPass.Name = "requirements - PASS"
$wsPass.Range("A1").Value = "Summary or requirements PASSED"
$styleSource.Copy()
$wsPass.Range("A1").PasteSpecial($xlPasteFormats)
$wsPass.PageSetup.CenterHeader = $headerCenter
$wsPass.PageSetup.LeftFooter = $footerLeft

# --- new sheet: requirements - FAIL ---------------------------------------
$wsFail = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsPass)
$wsFail.Name = "requirements - FAIL"
$wsFail.Range("A3").Value = "Summary or requirements FAIL"
$styleSource.Copy()
$wsFail.Range("A3").PasteSpecial($xlPasteFormats)
$wsFail.PageSetup.CenterHeader = $headerCenter
$wsFail.PageSetup.LeftFooter = $footerLeft

# --- new sheet: requirements - NULL ---------------------------------------
$wsNull = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsFail)
$wsNull.Name = "requirements - NULL"
$wsNull.Range("A5").Value = "Summary or requirements NULL"
$styleSource.Copy()
$wsNull.Range("A5").PasteSpecial($xlPasteFormats)
$wsNull.PageSetup.CenterHeader = $headerCenter
$wsNull.PageSetup.LeftFooter = $footerLeft

$excel.CutCopyMode = $false

# Keep "provenance" as the active/selected tab.
$ws1.Activate()
